$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stable donor cells (never themselves edited below) used purely to copy a
# fill/style onto a target cell via PasteSpecial(xlPasteFormats = -4122):
#   J2  -> style s=4 (dark-orange "room/adjacency" fill)
#   M3  -> style s=8 (grey "door#" fill)
$donor4 = $ws.Range("J2")
$donor8 = $ws.Range("M3")

# ---------------------------------------------------------------------------
# 1) U6: value changes to a brand new label "A*" (keeps its existing fill/style)
# ---------------------------------------------------------------------------
$ws.Range("U6").Value = "A*"

# ---------------------------------------------------------------------------
# 2) Row 10: K10 and L10 swap fill + text ("E" <-> "E#")
#    K10 was s=4/"E"  -> becomes s=8/"E#"
#    L10 was s=8/"E#" -> becomes s=4/"E"
# ---------------------------------------------------------------------------
$donor8.Copy() | Out-Null
$ws.Range("K10").PasteSpecial(-4122) | Out-Null
$ws.Range("K10").Value = "E#"

$donor4.Copy() | Out-Null
$ws.Range("L10").PasteSpecial(-4122) | Out-Null
$ws.Range("L10").Value = "E"

# ---------------------------------------------------------------------------
# 3) Row 20 / Row 21: D20 and F21 swap fill + text ("S" <-> "S#")
#    D20 was s=4/"S"  -> becomes s=8/"S#"
#    F21 was s=8/"S#" -> becomes s=4/"S"
# ---------------------------------------------------------------------------
$donor8.Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$ws.Range("D20").Value = "S#"

$donor4.Copy() | Out-Null
$ws.Range("F21").PasteSpecial(-4122) | Out-Null
$ws.Range("F21").Value = "S"

# ---------------------------------------------------------------------------
# 4) Row 21: T21 becomes "R#" (grey), U21/V21 become "R" (dark-orange)
# ---------------------------------------------------------------------------
$donor8.Copy() | Out-Null
$ws.Range("T21").PasteSpecial(-4122) | Out-Null
$ws.Range("T21").Value = "R#"

$donor4.Copy() | Out-Null
$ws.Range("U21").PasteSpecial(-4122) | Out-Null
$ws.Range("U21").Value = "R"

$donor4.Copy() | Out-Null
$ws.Range("V21").PasteSpecial(-4122) | Out-Null
$ws.Range("V21").Value = "R"

# ---------------------------------------------------------------------------
# 5) Row 22: T22, U22, V22 become "R" (dark-orange)
# ---------------------------------------------------------------------------
$donor4.Copy() | Out-Null
$ws.Range("T22").PasteSpecial(-4122) | Out-Null
$ws.Range("T22").Value = "R"

$donor4.Copy() | Out-Null
$ws.Range("U22").PasteSpecial(-4122) | Out-Null
$ws.Range("U22").Value = "R"

$donor4.Copy() | Out-Null
$ws.Range("V22").PasteSpecial(-4122) | Out-Null
$ws.Range("V22").Value = "R"

# ---------------------------------------------------------------------------
# 6) Row 26: X26 becomes "R" (dark-orange), was "R#" (grey)
# ---------------------------------------------------------------------------
$donor4.Copy() | Out-Null
$ws.Range("X26").PasteSpecial(-4122) | Out-Null
$ws.Range("X26").Value = "R"

# ---------------------------------------------------------------------------
# 7) Update the selection shown in the sheet view
# ---------------------------------------------------------------------------
$ws.Range("U7").Select() | Out-Null
